# Update gh-pages to output generated at 456a3b4
#
# This script:
#   1) Bumps the "想去人数" (F column) counters across the four sheets
#      (展览 / 演出 / 本地生活 / 全部类型) to their newly scraped values.
#   2) Appends a brand-new 展览 (exhibition) row for "北京·原神only"
#      at the end of the 展览 sheet (row 38).
#   3) Inserts the same new event as a row in the middle of 全部类型
#      (row 45), pushing the existing "Love Never Dies" row down to 46.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a date-looking string into a cell without Excel's
# automatic text->date coercion leaving a lingering number format on
# the cell (matches how the source data stores these as plain text).
# ---------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# =======================================================================
# 1) F-column ("想去人数") value bumps
# =======================================================================

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 3342
$ws1.Range("F6").Value  = 4837
$ws1.Range("F7").Value  = 469
$ws1.Range("F9").Value  = 176
$ws1.Range("F10").Value = 629
$ws1.Range("F11").Value = 283
$ws1.Range("F12").Value = 34
$ws1.Range("F13").Value = 14
$ws1.Range("F14").Value = 659
$ws1.Range("F15").Value = 289
$ws1.Range("F18").Value = 145
$ws1.Range("F19").Value = 344
$ws1.Range("F20").Value = 4758
$ws1.Range("F21").Value = 22
$ws1.Range("F22").Value = 35
$ws1.Range("F24").Value = 5901
$ws1.Range("F25").Value = 16
$ws1.Range("F26").Value = 1196
$ws1.Range("F27").Value = 240
$ws1.Range("F28").Value = 671
$ws1.Range("F29").Value = 4421
$ws1.Range("F32").Value = 126
$ws1.Range("F33").Value = 862
$ws1.Range("F34").Value = 74
$ws1.Range("F35").Value = 8
$ws1.Range("F36").Value = 790
$ws1.Range("F37").Value = 843

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 46

# --- 本地生活 (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 219
$ws3.Range("F4").Value = 40

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 219
$ws4.Range("F5").Value  = 40
$ws4.Range("F8").Value  = 3342
$ws4.Range("F10").Value = 4837
$ws4.Range("F11").Value = 469
$ws4.Range("F13").Value = 176
$ws4.Range("F14").Value = 629
$ws4.Range("F15").Value = 283
$ws4.Range("F16").Value = 34
$ws4.Range("F17").Value = 14
$ws4.Range("F18").Value = 659
$ws4.Range("F19").Value = 289
$ws4.Range("F23").Value = 145
$ws4.Range("F24").Value = 344
$ws4.Range("F25").Value = 4758
$ws4.Range("F26").Value = 22
$ws4.Range("F27").Value = 35
$ws4.Range("F29").Value = 5901
$ws4.Range("F30").Value = 16
$ws4.Range("F31").Value = 1196
$ws4.Range("F32").Value = 240
$ws4.Range("F33").Value = 671
$ws4.Range("F34").Value = 4421
$ws4.Range("F38").Value = 126
$ws4.Range("F39").Value = 862
$ws4.Range("F40").Value = 74
$ws4.Range("F41").Value = 8
$ws4.Range("F42").Value = 790
$ws4.Range("F43").Value = 843

# =======================================================================
# 2) 展览: append new row 38 ("北京·原神only")
# =======================================================================

$ws1.Rows.Item(38).Insert()
$ws1.Range("A37:I37").Copy()
$ws1.Range("A38:I38").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("A38").Value = 37
Set-TextValue $ws1.Range("B38") "2024-05-18"
$ws1.Range("C38").Value = "北京·原神only"
$ws1.Range("D38").Value = "北花园路1号 超级蜂巢"
$ws1.Range("E38").Value = "2024.05.18 10:00-05.19 17:00"
$ws1.Range("F38").Value = 1
$ws1.Range("G38").Value = 68
$ws1.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=81766"
$ws1.Range("I38").Value = "//i2.hdslb.com/bfs/openplatform/202402/Lfxwe5PO1707120983684.jpeg"

# =======================================================================
# 3) 全部类型: insert new row 45 ("北京·原神only"), pushing the old
#    row 45 (Love Never Dies) down to row 46.
# =======================================================================

$ws4.Rows.Item(45).Insert()
$ws4.Range("A44:I44").Copy()
$ws4.Range("A45:I45").PasteSpecial(-4122)   # xlPasteFormats

$ws4.Range("A45").Value = 44
Set-TextValue $ws4.Range("B45") "2024-05-18"
$ws4.Range("C45").Value = "北京·原神only"
$ws4.Range("D45").Value = "北花园路1号 超级蜂巢"
$ws4.Range("E45").Value = "2024.05.18 10:00-05.19 17:00"
$ws4.Range("F45").Value = 1
$ws4.Range("G45").Value = 68
$ws4.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=81766"
$ws4.Range("I45").Value = "//i2.hdslb.com/bfs/openplatform/202402/Lfxwe5PO1707120983684.jpeg"

# Row 46 (formerly row 45, "Love Never Dies") keeps its own data via the
# shift performed by Insert(); just fix up its index number in column A
# and bump its "want to go" counter (F) the same way F6 on 演出 was bumped
# (45 -> 46) for this same event.
$ws4.Range("A46").Value = 45
$ws4.Range("F46").Value = 46
